$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Pine"
$ws.Range("B3").Value = "Spruce"
$ws.Range("B4").Value = "Birch"
$ws.Range("B5").Value = "Other Deciduous"
